# Ccl20 -> Ackr4 NATMI TPM re-run: refresh computed metrics for the
# existing target clusters (ECs, FAPs, MuSCs) and append a new row for
# the "Resolving-Mac" target cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8243956666666667
$ws.Range("M2").Value = 0.011782
$ws.Range("N2").Value = 0.035346
$ws.Range("O2").Value = 0.0185453160301082
$ws.Range("P2").Value = 0.0190618769471875
$ws.Range("Q2").Value = 0.009713029744666668
$ws.Range("R2").Value = 0.08741726770200002
$ws.Range("S2").Value = 0.0185453160301082
$ws.Range("T2").Value = 0.0190618769471875

$ws.Range("G3").Value = 0.8243956666666667
$ws.Range("M3").Value = 0.450137
$ws.Range("O3").Value = 0.7085327552066554
$ws.Range("P3").Value = 0.7282682145116399
$ws.Range("Q3").Value = 0.3710909922063334
$ws.Range("S3").Value = 0.7085327552066554
$ws.Range("T3").Value = 0.7282682145116399

$ws.Range("G4").Value = 0.8243956666666667
$ws.Range("M4").Value = 0.051649
$ws.Range("N4").Value = 0.103298
$ws.Range("O4").Value = 0.08129749003896268
$ws.Range("P4").Value = 0.05570796596193557
$ws.Range("Q4").Value = 0.04257921178766667
$ws.Range("R4").Value = 0.255475270726
$ws.Range("S4").Value = 0.08129749003896268
$ws.Range("T4").Value = 0.05570796596193557

# New row for the Resolving-Mac target cluster
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl20"
$ws.Range("C5").Value = "Ackr4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8243956666666667
$ws.Range("H5").Value = 2.473187
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1217406666666667
$ws.Range("N5").Value = 0.365222
$ws.Range("O5").Value = 0.1916244387242736
$ws.Range("P5").Value = 0.1969619425792371
$ws.Range("Q5").Value = 0.1003624780571111
$ws.Range("R5").Value = 0.9032623025140001
$ws.Range("S5").Value = 0.1916244387242736
$ws.Range("T5").Value = 0.1969619425792371

